$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 28.02.2022 01:15"

# Row 5 (Makro) price update:
# New price goes to B5, previous B5 price shifts to C5 (as "Old Cena"),
# delta becomes a signed string in D5, and the check timestamp is written
# as a literal string in E5 (no date formatting).
$oldPrice = $ws.Range("B5").Value2

$ws.Range("B5").Value = 37.7
$ws.Range("C5").Value = $oldPrice

# D5/E5 must end up as plain text (inline string) cells with the
# workbook's default (General) style - i.e. no explicit NumberFormat
# should stick on the cell itself. Typing a leading "+" or a
# date-like string directly makes Excel auto-parse it back into a
# number, and pre-formatting the destination cell as Text leaves a
# lingering style on it. So: reset the destination's format back to
# the default first, then stage the literal text in a scratch cell
# that IS formatted as Text (so Excel keeps it as text when read back
# out of it), copy it, and paste-special *values only* into the
# destination - that carries over just the string while leaving the
# (already-default) destination format alone.
$scratch = $ws.Range("Z1")

$ws.Range("D5").ClearFormats()
$scratch.NumberFormat = "@"
$scratch.Value = "+0.4"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("E5").ClearFormats()
$scratch.NumberFormat = "@"
$scratch.Value = "2022-02-28 01:15:09"
$scratch.Copy()
$ws.Range("E5").PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
